# Apply "Latest stats from Matt" updates to the Poker - Year Figures workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update PersStatus for players who have become Inactive ---
# Every historical row belonging to pers_personid 354, 355 and 424 flips
# from "Active" to "Inactive".
$inactiveIds = @(354, 355, 424)
$lastRow = 199
for ($r = 2; $r -le $lastRow; $r++) {
    $personId = $ws.Cells.Item($r, 11).Value2
    if ($inactiveIds -contains $personId) {
        $ws.Cells.Item($r, 10).Value = "Inactive"
    }
}

# --- Refresh the 2024 year-to-date standings (rows 191-199) ---
$ws.Cells.Item(191, 1).Value = 2024
$ws.Cells.Item(191, 2).Value = "Richard"
$ws.Cells.Item(191, 3).Value = 1
$ws.Cells.Item(191, 4).Value = 55
$ws.Cells.Item(191, 5).Value = 0
$ws.Cells.Item(191, 6).Value = 55
$ws.Cells.Item(191, 7).Value = 174700
$ws.Cells.Item(191, 8).Value = 190
$ws.Cells.Item(191, 9).Value = 110
$ws.Cells.Item(191, 10).Value = "Active"
$ws.Cells.Item(191, 11).Value = 366

$ws.Cells.Item(192, 1).Value = 2024
$ws.Cells.Item(192, 2).Value = "Andy"
$ws.Cells.Item(192, 3).Value = 2
$ws.Cells.Item(192, 4).Value = 36
$ws.Cells.Item(192, 5).Value = 0
$ws.Cells.Item(192, 6).Value = 36
$ws.Cells.Item(192, 7).Value = 122600
$ws.Cells.Item(192, 8).Value = 130
$ws.Cells.Item(192, 9).Value = 50
$ws.Cells.Item(192, 10).Value = "Active"
$ws.Cells.Item(192, 11).Value = 349

$ws.Cells.Item(193, 1).Value = 2024
$ws.Cells.Item(193, 2).Value = "Mark"
$ws.Cells.Item(193, 3).Value = 3
$ws.Cells.Item(193, 4).Value = 36
$ws.Cells.Item(193, 5).Value = 0
$ws.Cells.Item(193, 6).Value = 36
$ws.Cells.Item(193, 7).Value = 111500
$ws.Cells.Item(193, 8).Value = 130
$ws.Cells.Item(193, 9).Value = 40
$ws.Cells.Item(193, 10).Value = "Active"
$ws.Cells.Item(193, 11).Value = 361

$ws.Cells.Item(194, 1).Value = 2024
$ws.Cells.Item(194, 2).Value = "Matt"
$ws.Cells.Item(194, 3).Value = 4
$ws.Cells.Item(194, 4).Value = 28
$ws.Cells.Item(194, 5).Value = 0
$ws.Cells.Item(194, 6).Value = 28
$ws.Cells.Item(194, 7).Value = 103550
$ws.Cells.Item(194, 8).Value = 40
$ws.Cells.Item(194, 9).Value = -50
$ws.Cells.Item(194, 10).Value = "Active"
$ws.Cells.Item(194, 11).Value = 362

$ws.Cells.Item(195, 1).Value = 2024
$ws.Cells.Item(195, 2).Value = "Anthony"
$ws.Cells.Item(195, 3).Value = 5
$ws.Cells.Item(195, 4).Value = 27
$ws.Cells.Item(195, 5).Value = 0
$ws.Cells.Item(195, 6).Value = 27
$ws.Cells.Item(195, 7).Value = 88050
$ws.Cells.Item(195, 8).Value = 40
$ws.Cells.Item(195, 9).Value = -20
$ws.Cells.Item(195, 10).Value = "Active"
$ws.Cells.Item(195, 11).Value = 350

$ws.Cells.Item(196, 1).Value = 2024
$ws.Cells.Item(196, 2).Value = "Pepe"
$ws.Cells.Item(196, 3).Value = 6
$ws.Cells.Item(196, 4).Value = 24
$ws.Cells.Item(196, 5).Value = 0
$ws.Cells.Item(196, 6).Value = 24
$ws.Cells.Item(196, 7).Value = 77350
$ws.Cells.Item(196, 8).Value = 40
$ws.Cells.Item(196, 9).Value = -40
$ws.Cells.Item(196, 10).Value = "Active"
$ws.Cells.Item(196, 11).Value = 364

$ws.Cells.Item(197, 1).Value = 2024
$ws.Cells.Item(197, 2).Value = "Jon"
$ws.Cells.Item(197, 3).Value = 7
$ws.Cells.Item(197, 4).Value = 22
$ws.Cells.Item(197, 5).Value = 0
$ws.Cells.Item(197, 6).Value = 22
$ws.Cells.Item(197, 7).Value = 69850
$ws.Cells.Item(197, 8).Value = 20
$ws.Cells.Item(197, 9).Value = -70
$ws.Cells.Item(197, 10).Value = "Active"
$ws.Cells.Item(197, 11).Value = 357

$ws.Cells.Item(198, 1).Value = 2024
$ws.Cells.Item(198, 2).Value = "Prashant"
$ws.Cells.Item(198, 3).Value = 8
$ws.Cells.Item(198, 4).Value = 20
$ws.Cells.Item(198, 5).Value = 0
$ws.Cells.Item(198, 6).Value = 20
$ws.Cells.Item(198, 7).Value = 76150
$ws.Cells.Item(198, 8).Value = 120
$ws.Cells.Item(198, 9).Value = 50
$ws.Cells.Item(198, 10).Value = "Active"
$ws.Cells.Item(198, 11).Value = 365

$ws.Cells.Item(199, 1).Value = 2024
$ws.Cells.Item(199, 2).Value = "Maisy"
$ws.Cells.Item(199, 3).Value = 9
$ws.Cells.Item(199, 4).Value = 14
$ws.Cells.Item(199, 5).Value = 0
$ws.Cells.Item(199, 6).Value = 14
$ws.Cells.Item(199, 7).Value = 55350
$ws.Cells.Item(199, 8).Value = 30
$ws.Cells.Item(199, 9).Value = -50
$ws.Cells.Item(199, 10).Value = "Active"
$ws.Cells.Item(199, 11).Value = 360

# --- New row 200: Alex finishes 10th for 2024 so far ---
$ws.Cells.Item(200, 1).Value = 2024
$ws.Cells.Item(200, 2).Value = "Alex"
$ws.Cells.Item(200, 3).Value = 10
$ws.Cells.Item(200, 4).Value = 4
$ws.Cells.Item(200, 5).Value = 0
$ws.Cells.Item(200, 6).Value = 4
$ws.Cells.Item(200, 7).Value = 17800
$ws.Cells.Item(200, 8).Value = 0
$ws.Cells.Item(200, 9).Value = -20
$ws.Cells.Item(200, 10).Value = "Active"
$ws.Cells.Item(200, 11).Value = 348

# --- Update the sheet dimension / selection to include the new row ---
$ws.Range("A1:K200").Select()
